$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2044887780548628
$ws.Range("C2").Value = 0.5461346633416458
$ws.Range("J2").Value = 0.002493765586034913
$ws.Range("P2").Value = 0.1471321695760598
$ws.Range("S2").Value = 0.09975062344139651
$ws.Range("B3").Value = 0.004504504504504504
$ws.Range("C3").Value = 0.004504504504504504
$ws.Range("J3").Value = 0.02702702702702703
$ws.Range("P3").Value = 0.7342342342342343
$ws.Range("S3").Value = 0.2297297297297297
$ws.Range("J4").Value = 0.08
$ws.Range("P4").Value = 0.72
$ws.Range("S4").Value = 0.2
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.06018518518518518
$ws.Range("D6").Value = 0.01388888888888889
$ws.Range("E6").Value = 0.004629629629629629
$ws.Range("F6").Value = 0.04629629629629629
$ws.Range("J6").Value = 0.3055555555555556
$ws.Range("O6").Value = 0.03240740740740741
$ws.Range("Q6").Value = 0.162037037037037
$ws.Range("R6").Value = 0.05555555555555555
$ws.Range("S6").Value = 0.3194444444444444
$ws.Range("B7").Value = 0.1525423728813559
$ws.Range("D7").Value = 0.00423728813559322
$ws.Range("F7").Value = 0.02966101694915254
$ws.Range("J7").Value = 0.1313559322033898
$ws.Range("O7").Value = 0.01694915254237288
$ws.Range("Q7").Value = 0.1949152542372881
$ws.Range("R7").Value = 0.06779661016949153
$ws.Range("S7").Value = 0.4025423728813559
$ws.Range("B8").Value = 0.1201814058956916
$ws.Range("D8").Value = 0.01360544217687075
$ws.Range("E8").Value = 0.006802721088435374
$ws.Range("F8").Value = 0.04308390022675737
$ws.Range("J8").Value = 0.1156462585034014
$ws.Range("Q8").Value = 0.2222222222222222
$ws.Range("R8").Value = 0.09070294784580499
$ws.Range("S8").Value = 0.3673469387755102
$ws.Range("B9").Value = 0.1388888888888889
$ws.Range("D9").Value = 0.02314814814814815
$ws.Range("F9").Value = 0.06018518518518518
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.01388888888888889
$ws.Range("Q9").Value = 0.1944444444444444
$ws.Range("R9").Value = 0.08333333333333333
$ws.Range("S9").Value = 0.3611111111111111
$ws.Range("B10").Value = 0.1277932960893855
$ws.Range("D10").Value = 0.0244413407821229
$ws.Range("E10").Value = 0.0006983240223463687
$ws.Range("F10").Value = 0.05446927374301676
$ws.Range("J10").Value = 0.1215083798882682
$ws.Range("O10").Value = 0.0111731843575419
$ws.Range("Q10").Value = 0.2129888268156425
$ws.Range("R10").Value = 0.08310055865921788
$ws.Range("S10").Value = 0.3638268156424581
$ws.Range("G11").Value = 0.1457725947521866
$ws.Range("J11").Value = 0.07580174927113703
$ws.Range("K11").Value = 0.2128279883381924
$ws.Range("L11").Value = 0.5393586005830904
$ws.Range("S11").Value = 0.02623906705539359
$ws.Range("G12").Value = 0.7323232323232324
$ws.Range("J12").Value = 0.1565656565656566
$ws.Range("K12").Value = 0.0202020202020202
$ws.Range("L12").Value = 0.06060606060606061
$ws.Range("S12").Value = 0.0303030303030303
$ws.Range("G13").Value = 0.8070175438596491
$ws.Range("J13").Value = 0.1754385964912281
$ws.Range("S13").Value = 0.01754385964912281
$ws.Range("F15").Value = 0.03153153153153153
$ws.Range("H15").Value = 0.1261261261261261
$ws.Range("I15").Value = 0.07207207207207207
$ws.Range("J15").Value = 0.3468468468468469
$ws.Range("K15").Value = 0.05855855855855856
$ws.Range("M15").Value = 0.009009009009009009
$ws.Range("N15").Value = 0.004504504504504504
$ws.Range("O15").Value = 0.08558558558558559
$ws.Range("S15").Value = 0.2657657657657658
$ws.Range("F16").Value = 0.04743083003952569
$ws.Range("H16").Value = 0.1343873517786561
$ws.Range("I16").Value = 0.1225296442687747
$ws.Range("J16").Value = 0.4268774703557312
$ws.Range("K16").Value = 0.08695652173913043
$ws.Range("M16").Value = 0.03557312252964427
$ws.Range("O16").Value = 0.03162055335968379
$ws.Range("S16").Value = 0.1146245059288538
$ws.Range("F17").Value = 0.01934235976789168
$ws.Range("H17").Value = 0.1411992263056093
$ws.Range("I17").Value = 0.0735009671179884
$ws.Range("J17").Value = 0.4526112185686654
$ws.Range("K17").Value = 0.09864603481624758
$ws.Range("M17").Value = 0.02707930367504836
$ws.Range("O17").Value = 0.05609284332688588
$ws.Range("S17").Value = 0.1315280464216634
$ws.Range("F18").Value = 0.01477832512315271
$ws.Range("H18").Value = 0.1625615763546798
$ws.Range("I18").Value = 0.08374384236453201
$ws.Range("J18").Value = 0.458128078817734
$ws.Range("K18").Value = 0.09359605911330049
$ws.Range("M18").Value = 0.01477832512315271
$ws.Range("O18").Value = 0.06896551724137931
$ws.Range("S18").Value = 0.103448275862069
$ws.Range("F19").Value = 0.02366863905325444
$ws.Range("H19").Value = 0.2026627218934911
$ws.Range("I19").Value = 0.08431952662721894
$ws.Range("J19").Value = 0.3764792899408284
$ws.Range("K19").Value = 0.1153846153846154
$ws.Range("M19").Value = 0.02144970414201183
$ws.Range("N19").Value = 0.0007396449704142012
$ws.Range("O19").Value = 0.06360946745562131
$ws.Range("S19").Value = 0.1116863905325444
